$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "charger mapping" table (currently A13:E18, with the charger-type
# labels in column F) down to A22:E27, dropping the column F labels entirely.

# 1) Copy the data + formatting of A13:E18 to the new location A22:E27.
$ws.Range("A13:E18").Copy($ws.Range("A22"))

# 2) Clear out the old block (A13:F18), including the column F charger-type
#    labels, so no residue (values/styles) is left behind at the old rows.
$ws.Range("A13:F18").Clear()

# 3) Restore the sheet selection to match the newly relocated table.
$ws.Range("A22:E27").Select()
